# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same B..H column layout) right after
#    itself, rename the copy to "2022-Q1" -> ends up positioned right
#    before "总计", and carries over the exact same cell styling.
# 2. Overwrite its values with the 2022-Q1 fund-holding data.
# 3. In "总计", grow the table by one row (copy the formatting of the
#    last existing data row down) and rewrite every data cell so the new
#    2022-Q1 summary row is on top and the rest shift down by one.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet as a copy of "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $sheets.Item("2021-Q4")
$q4.Copy([System.Reflection.Missing]::Value, $q4)
$newSheet = $sheets.Item($q4.Index + 1)
$newSheet.Name = "2022-Q1"

# Header row (kept identical to 2021-Q4, but D1 differs: 基金规模 not 基金金额)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Helper so text that looks numeric ("013776", "20.31", ...) is written
# verbatim as text instead of being coerced into a number.
function Set-TextCell($rng, $text) {
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2
$newSheet.Range("A2").Value = 0
Set-TextCell $newSheet.Range("B2") "013776"
Set-TextCell $newSheet.Range("C2") "中泰兴为价值精选混合A"
Set-TextCell $newSheet.Range("D2") "20.31"
Set-TextCell $newSheet.Range("E2") "85.34"
Set-TextCell $newSheet.Range("F2") "4.06"
Set-TextCell $newSheet.Range("G2") "0.8246"
$newSheet.Range("H2").Value = 10

# Row 3
$newSheet.Range("A3").Value = 1
Set-TextCell $newSheet.Range("B3") "010728"
Set-TextCell $newSheet.Range("C3") "中泰兴诚价值一年持有期混合A"
Set-TextCell $newSheet.Range("D3") "9.94"
Set-TextCell $newSheet.Range("E3") "88.08"
Set-TextCell $newSheet.Range("F3") "4.80"
Set-TextCell $newSheet.Range("G3") "0.4771"
$newSheet.Range("H3").Value = 10

# Row 4
$newSheet.Range("A4").Value = 2
Set-TextCell $newSheet.Range("B4") "013777"
Set-TextCell $newSheet.Range("C4") "中泰兴为价值精选混合C"
Set-TextCell $newSheet.Range("D4") "8.71"
Set-TextCell $newSheet.Range("E4") "85.34"
Set-TextCell $newSheet.Range("F4") "4.06"
Set-TextCell $newSheet.Range("G4") "0.3536"
$newSheet.Range("H4").Value = 10

# Row 5
$newSheet.Range("A5").Value = 3
Set-TextCell $newSheet.Range("B5") "010729"
Set-TextCell $newSheet.Range("C5") "中泰兴诚价值一年持有期混合C"
Set-TextCell $newSheet.Range("D5") "1.72"
Set-TextCell $newSheet.Range("E5") "88.08"
Set-TextCell $newSheet.Range("F5") "4.80"
Set-TextCell $newSheet.Range("G5") "0.0826"
$newSheet.Range("H5").Value = 10

# Row 6
$newSheet.Range("A6").Value = 4
Set-TextCell $newSheet.Range("B6") "007751"
Set-TextCell $newSheet.Range("C6") "景顺长城中证沪港深红利成长低波动指数A"
Set-TextCell $newSheet.Range("D6") "0.83"
Set-TextCell $newSheet.Range("E6") "91.29"
Set-TextCell $newSheet.Range("F6") "2.08"
Set-TextCell $newSheet.Range("G6") "0.0173"
$newSheet.Range("H6").Value = 10

# Row 7
$newSheet.Range("A7").Value = 5
Set-TextCell $newSheet.Range("B7") "007760"
Set-TextCell $newSheet.Range("C7") "景顺长城中证沪港深红利成长低波动指数C"
Set-TextCell $newSheet.Range("D7") "0.06"
Set-TextCell $newSheet.Range("E7") "91.29"
Set-TextCell $newSheet.Range("F7") "2.08"
Set-TextCell $newSheet.Range("G7") "0.0012"
$newSheet.Range("H7").Value = 10

# ---------------------------------------------------------------------
# Step 2: update the "总计" (totals) sheet with the new 2022-Q1 row
# ---------------------------------------------------------------------
$total = $sheets.Item("总计")

# Grow the table by one row, re-using the formatting of the last row.
$total.Range("A5:D5").Copy($total.Range("A6:D6"))

# Row 2 -> new 2022-Q1 summary (pushes everything else down by one row)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 1.76

# Row 3 -> was row 2 (2021-Q4)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 0.82

# Row 4 -> was row 3 (2021-Q3)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.72

# Row 5 -> was row 4 (2021-Q2)
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.65

# Row 6 -> was row 5 (2021-Q1)
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.03

Write-Host "2022-Q1 sheet added and 总计 refreshed"
